$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datadetails")

$ws.Range("B6").Value = "test1913"
$ws.Range("B7").Value = "Test1234@mailinator.com"
$ws.Range("B11").Value = "test1234"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "434575667"
